$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (A, C, D, F) keep their text formatting so that
# values such as dates and numeric-looking IDs are not auto-converted.
$ws.Range("A2:A34").NumberFormat = "@"
$ws.Range("C2:C34").NumberFormat = "@"
$ws.Range("D2:D34").NumberFormat = "@"
$ws.Range("F2:F34").NumberFormat = "@"

$ws.Range("A2").Value = "2025-07-07"
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = "OCRIM S A PRODUTOS ALIMENTICIOS"
$ws.Range("D2").Value = "53989115"
$ws.Range("E2").Value = 42648800
$ws.Range("F2").Value = "SACO DE LIXO 100L COMUM PACOTINHO C/ 5 UND CIDADE LIMPA/MAX LEVE"
$ws.Range("G2").Value = 1884
$ws.Range("H2").Value = 45.4
$ws.Range("I2").Value = 48.05

$ws.Range("A3").Value = "2025-07-07"
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = "ABC - AMAZONIA BRASIL CONCENTRADOS LTDA"
$ws.Range("D3").Value = "53975591"
$ws.Range("E3").Value = 27450023
$ws.Range("F3").Value = "LAMINA LARGA 18MM P/ ESTILETE C/10UND LEOARTE JOCAR"
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 3.17
$ws.Range("I3").Value = 2.2

$ws.Range("A4").Value = "2025-07-07"
$ws.Range("B4").Value = 350
$ws.Range("C4").Value = "SIND.DAS EMPRESAS DE TRANSP.DE PASSAGEIROS DO EST."
$ws.Range("D4").Value = "53958498"
$ws.Range("E4").Value = 30176476
$ws.Range("F4").Value = "ENVELOPE 24X34 OURO"
$ws.Range("G4").Value = 121
$ws.Range("H4").Value = 32
$ws.Range("I4").Value = 63.03

$ws.Range("A5").Value = "2025-07-07"
$ws.Range("B5").Value = 77
$ws.Range("C5").Value = "SIND.DAS EMPRESAS DE TRANSP.DE PASSAGEIROS DO EST."
$ws.Range("D5").Value = "53958498"
$ws.Range("E5").Value = 30240035
$ws.Range("F5").Value = "SACO DE LIXO 15L COMUM PACOTINHO C/20 UND"
$ws.Range("G5").Value = 195
$ws.Range("H5").Value = 13.97
$ws.Range("I5").Value = 15.79

$ws.Range("A6").Value = "2025-07-07"
$ws.Range("B6").Value = 600
$ws.Range("C6").Value = "POTENCIAL HUMANO RECRUTAMENTO E SELECAO LTDA"
$ws.Range("D6").Value = "53791039"
$ws.Range("E6").Value = 12775209
$ws.Range("F6").Value = "SACO DE LIXO 100L COMUM PACOTINHO C/5 UND PAPALIXO"
$ws.Range("G6").Value = 1076
$ws.Range("H6").Value = 37.38
$ws.Range("I6").Value = 61.82

$ws.Range("A7").Value = "2025-07-08"
$ws.Range("B7").Value = 15
$ws.Range("C7").Value = "CONDOMINIO RESIDENCIAL QUINTA DAS MARINAS"
$ws.Range("D7").Value = "54020891"
$ws.Range("E7").Value = 27437852
$ws.Range("F7").Value = "OCULOS DE PROTECAO SKY ESCURO CA 39878 DELTAPLUS"
$ws.Range("G7").Value = 28
$ws.Range("H7").Value = 4.07
$ws.Range("I7").Value = 3.25

$ws.Range("A8").Value = "2025-07-08"
$ws.Range("B8").Value = 100
$ws.Range("C8").Value = "MAYCA CONSTRUCOES E INSTALACOES ELETROMECANICAS EIRELI"
$ws.Range("D8").Value = "53473755"
$ws.Range("E8").Value = 12121233
$ws.Range("F8").Value = "SABAO EM PO ALA LAVANDA ROUPAS 400G"
$ws.Range("G8").Value = 102
$ws.Range("H8").Value = 16.52
$ws.Range("I8").Value = 25.02

$ws.Range("A9").Value = "2025-07-09"
$ws.Range("B9").Value = 200
$ws.Range("C9").Value = "V V REFEICOES LTDA"
$ws.Range("D9").Value = "54093449"
$ws.Range("E9").Value = 13811513
$ws.Range("F9").Value = "FIBRA DE LIMPEZA PESADA 102X260mm"
$ws.Range("G9").Value = 576
$ws.Range("H9").Value = 18.61
$ws.Range("I9").Value = 33.86

$ws.Range("A10").Value = "2025-07-09"
$ws.Range("B10").Value = 40
$ws.Range("C10").Value = "GREE ELECTRIC APPLIANCES DO BRASIL LTDA."
$ws.Range("D10").Value = "54093753"
$ws.Range("E10").Value = 12118255
$ws.Range("F10").Value = "LIMPA VIDROS GLOBO SAN 500ML"
$ws.Range("G10").Value = 905
$ws.Range("H10").Value = 8.01
$ws.Range("I10").Value = 10.58

$ws.Range("A11").Value = "2025-07-10"
$ws.Range("B11").Value = 270
$ws.Range("C11").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D11").Value = "54141318"
$ws.Range("E11").Value = 35118277
$ws.Range("F11").Value = "SABAO EM PO ABSOLUTO 400G"
$ws.Range("G11").Value = 711
$ws.Range("H11").Value = 18.27
$ws.Range("I11").Value = 33.82

$ws.Range("A12").Value = "2025-07-10"
$ws.Range("B12").Value = 95
$ws.Range("C12").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D12").Value = "54141318"
$ws.Range("E12").Value = 11939551
$ws.Range("F12").Value = "RODO COM CABO P 30CM"
$ws.Range("G12").Value = 50
$ws.Range("H12").Value = 6.4
$ws.Range("I12").Value = 15.76

$ws.Range("A13").Value = "2025-07-10"
$ws.Range("B13").Value = 95
$ws.Range("C13").Value = "SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA."
$ws.Range("D13").Value = "54142176"
$ws.Range("E13").Value = 19264853
$ws.Range("F13").Value = "PAPEL TOALHA INTERFOLHADO 1250 FLS NEWPAPER 100% CELULOSE"
$ws.Range("G13").Value = 142
$ws.Range("H13").Value = 8.29
$ws.Range("I13").Value = 13.53

$ws.Range("A14").Value = "2025-07-10"
$ws.Range("B14").Value = 30
$ws.Range("C14").Value = "CR OBRAS DA CONSTRUCAO LTDA"
$ws.Range("D14").Value = "54123729"
$ws.Range("E14").Value = 33278408
$ws.Range("F14").Value = "AROMATIZANTE LIMPADOR PERF CONC COALA ALGODAO 120ML"
$ws.Range("G14").Value = 36
$ws.Range("H14").Value = 5.92
$ws.Range("I14").Value = 5.81

$ws.Range("A15").Value = "2025-07-10"
$ws.Range("B15").Value = 35
$ws.Range("C15").Value = "SAT BRAS INDUSTRIA ELETRONICA DA AMAZONIA LTDA."
$ws.Range("D15").Value = "54142176"
$ws.Range("E15").Value = 11938367
$ws.Range("F15").Value = "PAPEL HIGIENICO 8X300 NEWPAPER 100% Celulose"
$ws.Range("G15").Value = 74
$ws.Range("H15").Value = 4.79
$ws.Range("I15").Value = 5.96

$ws.Range("A16").Value = "2025-07-10"
$ws.Range("B16").Value = 240
$ws.Range("C16").Value = "AMMAC INDUSTRIA E COMERCIO DE ALIMENTOS LTDA"
$ws.Range("D16").Value = "54153624"
$ws.Range("E16").Value = 32130390
$ws.Range("F16").Value = "ESPONJA MULTIUSO JEITOSA"
$ws.Range("G16").Value = 3388
$ws.Range("H16").Value = 21.06
$ws.Range("I16").Value = 50.05

$ws.Range("A17").Value = "2025-07-10"
$ws.Range("B17").Value = 112
$ws.Range("C17").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D17").Value = "54141318"
$ws.Range("E17").Value = 16537374
$ws.Range("F17").Value = "DESODORISADOR LADY AEROSSOL 360 ML LAVANDA"
$ws.Range("G17").Value = 2281
$ws.Range("H17").Value = 8.22
$ws.Range("I17").Value = 10.07

$ws.Range("A18").Value = "2025-07-10"
$ws.Range("B18").Value = 270
$ws.Range("C18").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D18").Value = "54141318"
$ws.Range("E18").Value = 14589837
$ws.Range("F18").Value = "FLANELA BRANCA TAM P 28X38CM"
$ws.Range("G18").Value = 283
$ws.Range("H18").Value = 15.43
$ws.Range("I18").Value = 28.28

$ws.Range("A19").Value = "2025-07-10"
$ws.Range("B19").Value = 113
$ws.Range("C19").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D19").Value = "54141318"
$ws.Range("E19").Value = 16871438
$ws.Range("F19").Value = "DESODORISADOR LADY AEROSSOL 360ML TALCO SUAVE CARINHO"
$ws.Range("G19").Value = 1481
$ws.Range("H19").Value = 9.23
$ws.Range("I19").Value = 11.4

$ws.Range("A20").Value = "2025-07-10"
$ws.Range("B20").Value = 89
$ws.Range("C20").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D20").Value = "54141318"
$ws.Range("E20").Value = 13995639
$ws.Range("F20").Value = "ESCOVA OVAL PLASTICA"
$ws.Range("G20").Value = 62
$ws.Range("H20").Value = 4.37
$ws.Range("I20").Value = 10.69

$ws.Range("A21").Value = "2025-07-10"
$ws.Range("B21").Value = 141
$ws.Range("C21").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D21").Value = "54141318"
$ws.Range("E21").Value = 12285275
$ws.Range("F21").Value = "LUSTRA MOVEIS BUTTERFLY 200ML AUDAX LAVANDA"
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 9.16
$ws.Range("I21").Value = 16.4

$ws.Range("A22").Value = "2025-07-10"
$ws.Range("B22").Value = 98
$ws.Range("C22").Value = "RODRIGUES INDUSTRIA E COMERCIO DE COLCHOES LTDA"
$ws.Range("D22").Value = "54141318"
$ws.Range("E22").Value = 11939681
$ws.Range("F22").Value = "VASSOURA VARRE CANTO COM CABO PLASTIFICADO"
$ws.Range("G22").Value = 85
$ws.Range("H22").Value = 5.85
$ws.Range("I22").Value = 11.19

$ws.Range("A23").Value = "2025-07-11"
$ws.Range("B23").Value = 310
$ws.Range("C23").Value = "JURUA ESTALEIROS E NAVEGACAO LTDA"
$ws.Range("D23").Value = "54127333"
$ws.Range("E23").Value = 12054191
$ws.Range("F23").Value = "PANO DE CHAO FLANELADO C REFORCADO ITATEX 42x62CM"
$ws.Range("G23").Value = 5764
$ws.Range("H23").Value = 20.03
$ws.Range("I23").Value = 39.31

$ws.Range("A24").Value = "2025-07-11"
$ws.Range("B24").Value = 154
$ws.Range("C24").Value = "JURUA ESTALEIROS E NAVEGACAO LTDA"
$ws.Range("D24").Value = "54124880"
$ws.Range("E24").Value = 11939672
$ws.Range("F24").Value = "VASSOURA PIACAVA 20 FUROS"
$ws.Range("G24").Value = 472
$ws.Range("H24").Value = 10.05
$ws.Range("I24").Value = 28.9

$ws.Range("A25").Value = "2025-07-14"
$ws.Range("B25").Value = 60
$ws.Range("C25").Value = "SAWEM DA AMAZONIA LTDA"
$ws.Range("D25").Value = "54285324"
$ws.Range("E25").Value = 27262762
$ws.Range("F25").Value = "MARCA TEXTO AMARELO UND JOCAR OFFICE"
$ws.Range("G25").Value = 111
$ws.Range("H25").Value = 7.25
$ws.Range("I25").Value = 8.13

$ws.Range("A26").Value = "2025-07-14"
$ws.Range("B26").Value = 250
$ws.Range("C26").Value = "ERAM ESTALEIRO RIO AMAZONAS LTDA"
$ws.Range("D26").Value = "54263882"
$ws.Range("E26").Value = 15011531
$ws.Range("F26").Value = "DETERGENTE LIMPOL COCO 500ML"
$ws.Range("G26").Value = 352
$ws.Range("H26").Value = 20.58
$ws.Range("I26").Value = 30.93

$ws.Range("A27").Value = "2025-07-15"
$ws.Range("B27").Value = 130
$ws.Range("C27").Value = "MUSASHI DA AMAZONIA LTDA"
$ws.Range("D27").Value = "54346779"
$ws.Range("E27").Value = 17171383
$ws.Range("F27").Value = "DETERGENTE DESENGRAX MAX PINE AUDAX 5L"
$ws.Range("G27").Value = -128
$ws.Range("H27").Value = 14.02
$ws.Range("I27").Value = 29.5

$ws.Range("A28").Value = "2025-07-15"
$ws.Range("B28").Value = 300
$ws.Range("C28").Value = "MUSASHI DA AMAZONIA LTDA"
$ws.Range("D28").Value = "54346779"
$ws.Range("E28").Value = 14795919
$ws.Range("F28").Value = "SACO DE LIXO 200L COMUM PACOTINHO C/5 UND SACOLMAX"
$ws.Range("G28").Value = 244
$ws.Range("H28").Value = 42.91
$ws.Range("I28").Value = 60.86

$ws.Range("A29").Value = "2025-07-16"
$ws.Range("B29").Value = 60
$ws.Range("C29").Value = "V V REFEICOES LTDA"
$ws.Range("D29").Value = "54396269"
$ws.Range("E29").Value = 42173656
$ws.Range("F29").Value = "COADOR DE CAFÉ G"
$ws.Range("G29").Value = 23
$ws.Range("H29").Value = 6.07
$ws.Range("I29").Value = 15.64

$ws.Range("A30").Value = "2025-07-16"
$ws.Range("B30").Value = 10
$ws.Range("C30").Value = "V V REFEICOES LTDA"
$ws.Range("D30").Value = "54396269"
$ws.Range("E30").Value = 17541022
$ws.Range("F30").Value = "COADOR DE CAFE INDUSTRIAL (MAIOR)"
$ws.Range("G30").Value = 29
$ws.Range("H30").Value = 2.08
$ws.Range("I30").Value = 2.06

$ws.Range("A31").Value = "2025-07-16"
$ws.Range("B31").Value = 10
$ws.Range("C31").Value = "V V REFEICOES LTDA"
$ws.Range("D31").Value = "54396269"
$ws.Range("E31").Value = 17125814
$ws.Range("F31").Value = "COADOR DE CAFE P"
$ws.Range("G31").Value = 21
$ws.Range("H31").Value = 2.19
$ws.Range("I31").Value = 1.98

$ws.Range("A32").Value = "2025-07-16"
$ws.Range("B32").Value = 70
$ws.Range("C32").Value = "TECHLOG - SERVICOS DE GESTAO E SISTEMAS INFORMATIZ"
$ws.Range("D32").Value = "54316914"
$ws.Range("E32").Value = 19264853
$ws.Range("F32").Value = "PAPEL TOALHA INTERFOLHADO 1250 FLS NEWPAPER 100% CELULOSE"
$ws.Range("G32").Value = 142
$ws.Range("H32").Value = 8.29
$ws.Range("I32").Value = 13.53

$ws.Range("A33").Value = "2025-07-16"
$ws.Range("B33").Value = 10
$ws.Range("C33").Value = "CONDOMINIO DO TVLANDIA MALL"
$ws.Range("D33").Value = "54119372"
$ws.Range("E33").Value = 28133466
$ws.Range("F33").Value = "DESINFETANTE CONCENTRADO 5L AUDAX MAX 1:200 - LAVANDA"
$ws.Range("G33").Value = 5
$ws.Range("H33").Value = 2.34
$ws.Range("I33").Value = 1.91

$ws.Range("A34").Value = "2025-07-17"
$ws.Range("B34").Value = 100
$ws.Range("C34").Value = "TEL TELECOMUNICACOES LTDA."
$ws.Range("D34").Value = "54443314"
$ws.Range("E34").Value = 11936640
$ws.Range("F34").Value = "LIMPADOR VEJA MULTIUSO GOLD 500ML"
$ws.Range("G34").Value = 3426
$ws.Range("H34").Value = 10.58
$ws.Range("I34").Value = 14.64

# Remove the now-obsolete last row (old row 35) so the table has 34 rows total
$ws.Rows(35).Delete()

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()